# Modified DSL for EB
# Clear out the Pass/Fail markers that used to live in column J (rows 2-5)
# of the TestCases sheet, and move the active selection to E2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

$ws.Range("J2").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("J5").ClearContents()

$ws.Activate()
$ws.Range("E2").Select()
